$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the next day's case record to the dataset (row 28)
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 43918
$ws.Range("B28").Value = 1203
$ws.Range("C28").Value = 99
$ws.Range("D28").Value = 10
$ws.Range("E28").Value = 89
$ws.Range("F28").Value = 0

# Match the saved selection state from the edit
$ws.Range("F28").Select()
